$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '243.08'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '23.71'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.238'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05757'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.405'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8054'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8884'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1374'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07079'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03041'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09310'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.807'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001536'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04718'
$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0006020'
$ws.Range("E18").Value = '17OneONEWorstin24h'
$ws.Range("B19").Value = 'TigerCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.006182'
$ws.Range("E19").Value = '18TigerCashTCH'
$ws.Range("B20").Value = 'BitKan'
$ws.Range("C20").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.001249'
$ws.Range("E20").Value = '19BitKanKAN'
$ws.Range("B21").Value = 'HotbitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.004061'
$ws.Range("E21").Value = '20HotbitTokenHTB'
$ws.Range("B22").Value = 'NitroEx'
$ws.Range("C22").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.00008698'
$ws.Range("E22").Value = '21NitroExNTX'
$ws.Range("B23").Value = 'LEO'
$ws.Range("C23").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.543'
$ws.Range("E23").Value = '22LEOLEO'
$ws.Range("B24").Value = 'BTSEToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.140'
$ws.Range("E24").Value = '23BTSETokenBTSE'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1310'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03717'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006239'
$ws.Range("E41").Value = '40KickTokenKICK'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1046'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002491'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.007131'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005327'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5350'
$ws.Range("E47").Value = '46CoinbaseStockTokenCOIN'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002874'
$ws.Range("E48").Value = '47BOLOBOLOBestin24h'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002100'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0002000'
